{"js": "// The edit applies strikethrough formatting to the five short paragraphs\n// near the top of the document (the \"Color scheme\" / font notes block),\n// i.e. everything after the \"Jenn info\" title and up to (and including)\n// \"Poppins Extra Light\" -- but NOT the \"Paypal ...\" paragraph that follows.\n// Setting Font.strikeThrough on a paragraph writes <w:strike/> into both\n// the paragraph mark's rPr (w:pPr/w:rPr) and every run's rPr, matching the\n// target diff exactly.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The block of paragraphs to strike through, identified by their exact text\n// (robust against accidental reordering / insertion elsewhere in the doc).\nconst targetTexts = [\n  \"Color scheme: Triad RGB codes 64,189, 26 & 26,64,189 & 189,26,64\",\n  \"I'm not sure how many we need? Here are 3 that may pair well if needed. Fonts\",\n  \"Reclame Script I could not find. But I read this combined well with the next two. I looked up sister fonts as well--bakerie & brush script and could not find. It's like a street style font. If you can't find it either we can go with this one Euphoria script\",\n  \"Poppins Semi Bold\",\n  \"Poppins Extra Light\",\n];\n\nfor (const paragraph of paragraphs.items) {\n  if (targetTexts.includes(paragraph.text)) {\n    paragraph.font.strikeThrough = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The edit applies strikethrough formatting to the five short paragraphs\n# near the top of the document (the \"Color scheme\" / font notes block),\n# i.e. everything after the \"Jenn info\" title and up to (and including)\n# \"Poppins Extra Light\" -- but NOT the \"Paypal ...\" paragraph that follows.\n# Setting Range.Font.StrikeThrough on a paragraph's range writes <w:strike/>\n# into both the paragraph mark's rPr (w:pPr/w:rPr) and every run's rPr,\n# matching the target diff exactly.\n\n$d = $word.ActiveDocument\n\n$targetTexts = @(\n  \"Color scheme: Triad RGB codes 64,189, 26 & 26,64,189 & 189,26,64\",\n  \"I'm not sure how many we need? Here are 3 that may pair well if needed. Fonts\",\n  \"Reclame Script I could not find. But I read this combined well with the next two. I looked up sister fonts as well--bakerie & brush script and could not find. It's like a street style font. If you can't find it either we can go with this one Euphoria script\",\n  \"Poppins Semi Bold\",\n  \"Poppins Extra Light\"\n)\n\nforeach ($p in $d.Paragraphs) {\n  # Paragraph.Range.Text includes the trailing paragraph-mark character\n  # (CR, char 13) -- strip it before comparing against the target text.\n  $text = $p.Range.Text.TrimEnd([char]13, [char]10)\n  if ($targetTexts -contains $text) {\n    $p.Range.Font.StrikeThrough = 1\n  }\n}\n"}
